$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A36").Value = 973
$ws.Range("B36").Value = "K Closest Points to Origin"
